# The deck ships two theme parts: ppt/theme/theme1.xml ("Office Theme",
# used by the Notes Master) and ppt/theme/theme2.xml ("Integral", the
# theme actually applied to the Slide Master / all slides). The edit
# swaps the two themes' content, so the slides switch from the
# "Integral" palette to the standard "Office Theme" palette.
#
# The PowerPoint object model only exposes a single Theme for editing
# (reached via the Slide Master / Design), which is backed by the part
# driving the visible slides (theme2.xml). We reapply the target
# "Office Theme" color scheme onto it, one swatch at a time, via
# ThemeColorScheme -- the COM-visible equivalent of editing <a:clrScheme>.

$p = $ppt.ActivePresentation
$design = $p.Designs.Item(1)
$theme = $design.SlideMaster.Theme
$colors = $theme.ThemeColorScheme

# Office Theme color scheme (dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink)
$colors.Item(1).RGB  = 0x00 + 0x00*256 + 0x00*65536    # dk1      000000
$colors.Item(2).RGB  = 0xFF + 0xFF*256 + 0xFF*65536    # lt1      FFFFFF
$colors.Item(3).RGB  = 0x44 + 0x54*256 + 0x6A*65536    # dk2      44546A
$colors.Item(4).RGB  = 0xE7 + 0xE6*256 + 0xE6*65536    # lt2      E7E6E6
$colors.Item(5).RGB  = 0x5B + 0x9B*256 + 0xD5*65536    # accent1  5B9BD5
$colors.Item(6).RGB  = 0xED + 0x7D*256 + 0x31*65536    # accent2  ED7D31
$colors.Item(7).RGB  = 0xA5 + 0xA5*256 + 0xA5*65536    # accent3  A5A5A5
$colors.Item(8).RGB  = 0xFF + 0xC0*256 + 0x00*65536    # accent4  FFC000
$colors.Item(9).RGB  = 0x44 + 0x72*256 + 0xC4*65536    # accent5  4472C4
$colors.Item(10).RGB = 0x70 + 0xAD*256 + 0x47*65536    # accent6  70AD47
$colors.Item(11).RGB = 0x05 + 0x63*256 + 0xC1*65536    # hlink    0563C1
$colors.Item(12).RGB = 0x95 + 0x4F*256 + 0x72*65536    # folHlink 954F72
